$d = $word.ActiveDocument

# Locate (1-based) the paragraph that holds the "Ver no Jupiter ..." line and
# the one that holds the "(c) <year> . Contact: ..." footer line that follows
# it. These two paragraphs -- together with the blank paragraph that
# precedes "Ver no Jupiter" -- are the ones being removed by this edit. The
# blank paragraph that originally followed the copyright line (the one right
# before the page-break paragraph) is left untouched.
$jupiterIndex = -1
$copyrightIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $jupiterIndex = $i
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $copyrightIndex = $i
    }
    $i = $i + 1
}

if ($jupiterIndex -gt 0) {
    # Work out how many paragraphs must disappear: the optional blank one
    # right above "Ver no Jupiter", plus every paragraph from "Ver no
    # Jupiter" through the copyright line (inclusive).
    $firstToRemove = $jupiterIndex
    if ($jupiterIndex -gt 1) {
        $prevText = $d.Paragraphs.Item($jupiterIndex - 1).Range.Text
        # A paragraph's Range.Text includes its trailing paragraph mark, so
        # an "empty" paragraph's text is just that single mark character.
        if ($prevText.Length -le 1) {
            $firstToRemove = $jupiterIndex - 1
        }
    }

    if ($copyrightIndex -ge $jupiterIndex) {
        $lastToRemove = $copyrightIndex
    } else {
        $lastToRemove = $jupiterIndex
    }

    $removeCount = $lastToRemove - $firstToRemove + 1

    for ($n = 0; $n -lt $removeCount; $n++) {
        $d.Paragraphs.Item($firstToRemove).Range.Delete()
    }
}
